$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Raw input value corrections (no formulas involved) ---
$ws.Range("D12").Value = 1398800851.3300049
$ws.Range("E12").Value = 1361974149

$ws.Range("D13").Value = 338965688.90999979
$ws.Range("E13").Value = 332114255.60000002

$ws.Range("D14").Value = -45752811.059999987
$ws.Range("E14").Value = 537985.22

$ws.Range("D16").Value = -50601311.959999993
$ws.Range("E16").Value = -60413972.810000002

$ws.Range("D19").Value = -383099999.99999988
$ws.Range("E19").Value = -384700000

$ws.Range("D26").Value = 1009991810.1331247
$ws.Range("E26").Value = 1018613404

# --- D18 and D21 were plain numbers; now become SUM formulas like their neighbors ---
$ws.Range("D18").Formula = "=SUM(D12:D17)"
$ws.Range("D21").Formula = "=SUM(D18:D20)"

# Recalculate the whole workbook so dependent formula cells (E18, D21/E21, D23/E23,
# D25/E25, D28/E28/F28, C29, etc.) pick up the new cached values.
$excel.CalculateFullRebuild()
